$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held a small 2-row "Verhoudingen" lookup table:
#   A1 = "op ware grootte"
#   A2 = "Verhoudingen"
# This upload fixes the capitalisation of A1 and appends two more rows
# ("Percentage" and "verhoudingstabel") turning it into a 4-row table.
#
# Note: the write order below is deliberate. Writing "Percentage" into A1
# first reuses/edits the shared-string slot that used to hold
# "op ware grootte", then A3 is given its own "Percentage" string entry,
# and only afterwards is A1 corrected to its real final value
# ("Op ware grootte"). That reproduces the exact shared-string insertion
# order of the authored workbook (Verhoudingen, Percentage,
# Op ware grootte, verhoudingstabel).
$ws.Range("A1").Value = "Percentage"
$ws.Range("A2").Value = "Verhoudingen"
$ws.Range("A3").Value = "Percentage"
$ws.Range("A1").Value = "Op ware grootte"
$ws.Range("A4").Value = "verhoudingstabel"

# Final selection lands on the next empty row, same as in the saved file.
$ws.Range("A5").Select()
